$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text representation (avoid Excel auto-converting
# numeric-looking strings like "1.00" or "0.100" into numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.587.50'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.913.49'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.47'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.96%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.95'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +8.73%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.59%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.100'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.55'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +8.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.815'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.191.64'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.11'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.915.02'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.633.17'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.34'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0859'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '249.80'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.32'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.20'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.64'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.20'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.87'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.81'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.69'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.18%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.59'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.02%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.95'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.08%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0881'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +20.72%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.53'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.877'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.57'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +49.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.03'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.57'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +10.51%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.35'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.82'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +17.84%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.10'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.343.22'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.64%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0813'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.80'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.39'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.46'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.28%  '
